# Update column F (dSF) values on the active sheet to reflect the
# repulled/pushed data mentioned in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 0
    3  = -3
    6  = 2
    8  = 0
    9  = -1
    11 = -3
    12 = -7
    14 = -4
    17 = -10
    20 = -4
    23 = 8
    25 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
